$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.88321066666667
$ws.Range("H2").Value = 107.649632
$ws.Range("I2").Value = 0.08317795499144418
$ws.Range("J2").Value = 0.08448843719082051
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 36.23704304598045
$ws.Range("R2").Value = 326.133387413824
$ws.Range("S2").Value = 0.001125741759813923
$ws.Range("T2").Value = 0.001303251946201785
$ws.Range("G3").Value = 35.88321066666667
$ws.Range("H3").Value = 107.649632
$ws.Range("I3").Value = 0.08317795499144418
$ws.Range("J3").Value = 0.08448843719082051
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 1656.145117060914
$ws.Range("R3").Value = 14905.30605354822
$ws.Range("S3").Value = 0.05144988558315038
$ws.Range("T3").Value = 0.05956265096640204
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 35.88321066666667
$ws.Range("H4").Value = 107.649632
$ws.Range("I4").Value = 0.08317795499144418
$ws.Range("J4").Value = 0.08448843719082051
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.3342760295004445
$ws.Range("R4").Value = 3.008484265504
$ws.Range("S4").Value = 0.00001038463555748606
$ws.Range("T4").Value = 0.00001202211464832485
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 35.88321066666667
$ws.Range("H5").Value = 107.649632
$ws.Range("I5").Value = 0.08317795499144418
$ws.Range("J5").Value = 0.08448843719082051
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 984.7387699312908
$ws.Range("R5").Value = 5908.432619587745
$ws.Range("S5").Value = 0.03059194301292239
$ws.Range("T5").Value = 0.02361051216356835
$ws.Range("I6").Value = 0.03522729558434242
$ws.Range("J6").Value = 0.03578230735158529
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 15.34701143607822
$ws.Range("R6").Value = 138.123102924704
$ws.Range("S6").Value = 0.0004767710113657165
$ws.Range("T6").Value = 0.0005519496305775034
$ws.Range("I7").Value = 0.03522729558434242
$ws.Range("J7").Value = 0.03578230735158529
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.02178991209154717
$ws.Range("T7").Value = 0.02522580786695575
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.03522729558434242
$ws.Range("J8").Value = 0.03578230735158529
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 0.1415716519982222
$ws.Range("R8").Value = 1.274144867984
$ws.Range("S8").Value = 0.000004398071897257668
$ws.Range("T8").Value = 0.000005091572476252331
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.03522729558434242
$ws.Range("J9").Value = 0.03578230735158529
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 417.0538182298373
$ws.Range("R9").Value = 2502.322909379024
$ws.Range("S9").Value = 0.01295621440953227
$ws.Range("T9").Value = 0.009999458281575788
$ws.Range("G10").Value = 177.70077
$ws.Range("H10").Value = 533.10231
$ws.Range("I10").Value = 0.4119137160358794
$ws.Range("J10").Value = 0.4184034835782469
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 179.45301805938
$ws.Range("R10").Value = 1615.07716253442
$ws.Range("S10").Value = 0.005574896276656734
$ws.Range("T10").Value = 0.006453961895867581
$ws.Range("G11").Value = 177.70077
$ws.Range("H11").Value = 533.10231
$ws.Range("I11").Value = 0.4119137160358794
$ws.Range("J11").Value = 0.4184034835782469
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 8201.558808862379
$ws.Range("R11").Value = 73814.02927976141
$ws.Range("S11").Value = 0.2547900289488511
$ws.Range("T11").Value = 0.2949660507888467
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 177.70077
$ws.Range("H12").Value = 533.10231
$ws.Range("I12").Value = 0.4119137160358794
$ws.Range("J12").Value = 0.4184034835782469
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 1.65540113973
$ws.Range("R12").Value = 14.89861025757
$ws.Range("S12").Value = 0.00005142677314683206
$ws.Range("T12").Value = 0.00005953589409489866
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 177.70077
$ws.Range("H13").Value = 533.10231
$ws.Range("I13").Value = 0.4119137160358794
$ws.Range("J13").Value = 0.4184034835782469
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 4876.621528970295
$ws.Range("R13").Value = 29259.72917382177
$ws.Range("S13").Value = 0.1514973640372248
$ws.Range("T13").Value = 0.1169239349994377
$ws.Range("G14").Value = 20.074196
$ws.Range("H14").Value = 40.148392
$ws.Range("I14").Value = 0.04653236263856699
$ws.Range("J14").Value = 0.0315103250497358
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 20.27214095535734
$ws.Range("R14").Value = 121.632845732144
$ws.Range("S14").Value = 0.0006297753270133692
$ws.Range("T14").Value = 0.0004860534034233595
$ws.Range("G15").Value = 20.074196
$ws.Range("H15").Value = 40.148392
$ws.Range("I15").Value = 0.04653236263856699
$ws.Range("J15").Value = 0.0315103250497358
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 926.4996377597573
$ws.Range("R15").Value = 5558.997826558543
$ws.Range("S15").Value = 0.02878268327123687
$ws.Range("T15").Value = 0.02221414616223015
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 20.074196
$ws.Range("H16").Value = 40.148392
$ws.Range("I16").Value = 0.04653236263856699
$ws.Range("J16").Value = 0.0315103250497358
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.1870045185373333
$ws.Range("R16").Value = 1.122027111224
$ws.Range("S16").Value = 0.000005809491561556225
$ws.Range("T16").Value = 0.000004483699224999563
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 20.074196
$ws.Range("H17").Value = 40.148392
$ws.Range("I17").Value = 0.04653236263856699
$ws.Range("J17").Value = 0.0315103250497358
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 550.8938221841661
$ws.Range("R17").Value = 2203.575288736664
$ws.Range("S17").Value = 0.0171140945487552
$ws.Range("T17").Value = 0.008805641784857288
$ws.Range("G18").Value = 182.547562
$ws.Range("H18").Value = 547.642686
$ws.Range("I18").Value = 0.423148670749767
$ws.Range("J18").Value = 0.4298154468296114
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 184.3476026596947
$ws.Range("R18").Value = 1659.128423937252
$ws.Range("S18").Value = 0.005726951682350979
$ws.Range("T18").Value = 0.006629993833631249
$ws.Range("G19").Value = 182.547562
$ws.Range("H19").Value = 547.642686
$ws.Range("I19").Value = 0.423148670749767
$ws.Range("J19").Value = 0.4298154468296114
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 8425.256486831495
$ws.Range("R19").Value = 75827.30838148344
$ws.Range("S19").Value = 0.2617394320042742
$ws.Range("T19").Value = 0.3030112556308683
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 182.547562
$ws.Range("H20").Value = 547.642686
$ws.Range("I20").Value = 0.423148670749767
$ws.Range("J20").Value = 0.4298154468296114
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 1.700552238404667
$ws.Range("R20").Value = 15.304970145642
$ws.Range("S20").Value = 0.0000528294393968088
$ws.Range("T20").Value = 0.00006115973677837158
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 182.547562
$ws.Range("H21").Value = 547.642686
$ws.Range("I21").Value = 0.423148670749767
$ws.Range("J21").Value = 0.4298154468296114
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 5009.631477175027
$ws.Range("R21").Value = 30057.78886305017
$ws.Range("S21").Value = 0.155629457623745
$ws.Range("T21").Value = 0.1201130376283334
